# Fill rows 3 and 4 with the same flight data as row 2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 3, 4) {
    $ws.Cells.Item($r, 1).Value = "VN0012"
    $ws.Cells.Item($r, 2).Value = "VN4567"
    $ws.Cells.Item($r, 3).Value = 45635.375
    $ws.Cells.Item($r, 4).Value = 90
    $ws.Cells.Item($r, 5).Value = "Tan Son Nhat International Airport"
    $ws.Cells.Item($r, 6).Value = "Noi Bai International Airport"
    $ws.Cells.Item($r, 7).Value = 1490000
    $ws.Cells.Item($r, 8).Value = 2000000
    $ws.Cells.Item($r, 9).Value = 50000000
}

# Move the active selection to A11 as recorded in the saved view state
$ws.Range("A11").Select()
